$wb = $excel.ActiveWorkbook

# The "Events-Irrigation" sheet carries the ISA annotation table
# ("annotationTable"). Two of its columns are being relabeled from
# "... [Source Name]" to "... [Sample Name]" to better reflect that the
# input/output material is a Sample, not a raw Source.
$ws = $wb.Worksheets.Item("Events-Irrigation")

# Renaming the header cell of a ListObject column renames the column
# itself (ListColumns(i).Name follows the header text).
$ws.Range("A1").Value = "Input [Sample Name]"
$ws.Range("V1").Value = "Output [Sample Name]"
